$p = $ppt.ActivePresentation

# Unicode quote characters used throughout the deck's reference list.
$lq = [char]0x201C   # “
$rq = [char]0x201D   # ”

# -----------------------------------------------------------------
# 1) Slide 29 ("ROSA [17]" -> "ROSA [17], [30], [31]")
# -----------------------------------------------------------------
$s29 = $p.Slides.Item(29)
$sh29 = $s29.Shapes.Item(3)
$tr29 = $sh29.TextFrame.TextRange
$full29 = $tr29.Text
$needle29 = "ROSA [17]"
$idx29 = $full29.IndexOf($needle29)
if ($idx29 -ge 0) {
    $sub29 = $tr29.Characters($idx29 + 1, $needle29.Length)
    $sub29.Text = "ROSA [17], [30], [31]"
}

# -----------------------------------------------------------------
# 2) Slide 53 ("References") - italicize a title, then add [31] and
#    [32] reference paragraphs.
# -----------------------------------------------------------------
$s53 = $p.Slides.Item(53)
$sh53 = $s53.Shapes.Item(2)
$tr53 = $sh53.TextFrame.TextRange

# 2a) Italicize "U.S. Bureau of Labor Statistics" inside the [30] entry.
$full53 = $tr53.Text
$needleA = "U.S. Bureau of Labor Statistics"
$idxA = $full53.IndexOf($needleA)
if ($idxA -ge 0) {
    $subA = $tr53.Characters($idxA + 1, $needleA.Length)
    $subA.Font.Italic = $true
}

# 2b) Build the text for the two new reference paragraphs, plus one new
#     blank trailing paragraph, and insert them right after paragraph 1
#     (the [30] entry) so the pre-existing blank paragraph stays last.
$para1 = $tr53.Paragraphs(1, 1)

$oUml = [char]0x00F6   # o
$gBreve = [char]0x011F # g
$dotlessI = [char]0x0131  # i
$cCedil = [char]0x00E7 # c

$gosterge = "G" + $oUml + "sterge"
$niteligindeki = "Niteli" + $gBreve + "indeki"
$bankasi = "Bankas" + $dotlessI
$kurlari = "Kurlar" + $dotlessI
$kac = "ka" + $cCedil

$p31 = " [31] " + $lq + $gosterge + " " + $niteligindeki + " Merkez " + $bankasi + " " + $kurlari + "," + $rq + " TCMB. [Online]. Available: https://www.tcmb.gov.tr/kurlar/kurlar_tr.html https://www.bls.gov/ooh/architecture-and-engineering/mechanical-engineers.html. [Accessed: Mar. 23, 2023]."
$p32 = " [32] " + $lq + "1 kWh Elektrik " + $kac + " TL ?," + $rq + " GazElektrik. [Online]. Available: https://gazelektrik.com/faydali-bilgiler/1-kwh-elektrik-kac-tl. [Accessed: Mar. 23, 2023]."

$insertText = "`r" + $p31 + "`r" + $p32 + "`r"
$para1.InsertAfter($insertText) | Out-Null

$fullAfter = $tr53.Text

# 2c) Italicize "TCMB" within the new [31] paragraph.
$idxTcmb = $fullAfter.IndexOf("TCMB")
if ($idxTcmb -ge 0) {
    $subTcmb = $tr53.Characters($idxTcmb + 1, 4)
    $subTcmb.Font.Italic = $true
}

# 2d) Italicize "GazElektrik" within the new [32] paragraph.
$idxGaz = $fullAfter.IndexOf("GazElektrik")
if ($idxGaz -ge 0) {
    $subGaz = $tr53.Characters($idxGaz + 1, "GazElektrik".Length)
    $subGaz.Font.Italic = $true
}

# 2e) Add the hyperlink to the gazelektrik URL text in the [32] paragraph.
$url = "https://gazelektrik.com/faydali-bilgiler/1-kwh-elektrik-kac-tl"
$idxUrl = $fullAfter.IndexOf($url)
if ($idxUrl -ge 0) {
    $subUrl = $tr53.Characters($idxUrl + 1, $url.Length)
    $subUrl.ActionSettings(1).Hyperlink.Address = $url
}
